$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically rotates the species-identifying data of rows 5, 6 and 7:
#   old row 7 -> new row 5
#   old row 5 -> new row 6
#   old row 6 -> new row 7
# (all other columns in these rows stay identical, so only the columns
# A, B, E, F, G, H, L, Q, R need to be rewritten)

$row5 = @{
    A = 111413026
    B = 89007
    E = 1068
    F = "Liten stinksvamp"
    G = "Mutinus caninus"
    H = "(Schaeff.:Pers.) Fr."
    Q = 492828.1855948549
    R = 6227996.970613244
}

$row6 = @{
    A = 111412885
    B = 90536
    E = 1428
    F = "Rutkremla"
    G = "Russula virescens"
    H = "(Schaeff.) Fr."
    Q = 493001.1390786725
    R = 6227751.92766118
}

$row7 = @{
    A = 111412863
    B = 98535
    E = 222498
    F = "Blåsippa"
    G = "Hepatica nobilis"
    H = "Schreb."
    Q = 493001.1390786725
    R = 6227751.92766118
}

# New row 5 gets the old row 7 values
$ws.Cells.Item(5, 1).Value  = $row7.A
$ws.Cells.Item(5, 2).Value  = $row7.B
$ws.Cells.Item(5, 5).Value  = $row7.E
$ws.Cells.Item(5, 6).Value  = $row7.F
$ws.Cells.Item(5, 7).Value  = $row7.G
$ws.Cells.Item(5, 8).Value  = $row7.H
$ws.Cells.Item(5, 17).Value = $row7.Q
$ws.Cells.Item(5, 18).Value = $row7.R

# New row 6 gets the old row 5 values
$ws.Cells.Item(6, 1).Value  = $row5.A
$ws.Cells.Item(6, 2).Value  = $row5.B
$ws.Cells.Item(6, 5).Value  = $row5.E
$ws.Cells.Item(6, 6).Value  = $row5.F
$ws.Cells.Item(6, 7).Value  = $row5.G
$ws.Cells.Item(6, 8).Value  = $row5.H
$ws.Cells.Item(6, 17).Value = $row5.Q
$ws.Cells.Item(6, 18).Value = $row5.R

# New row 7 gets the old row 6 values
$ws.Cells.Item(7, 1).Value  = $row6.A
$ws.Cells.Item(7, 2).Value  = $row6.B
$ws.Cells.Item(7, 5).Value  = $row6.E
$ws.Cells.Item(7, 6).Value  = $row6.F
$ws.Cells.Item(7, 7).Value  = $row6.G
$ws.Cells.Item(7, 8).Value  = $row6.H
$ws.Cells.Item(7, 17).Value = $row6.Q
$ws.Cells.Item(7, 18).Value = $row6.R

# Column L ("Kön") is an empty marker cell that used to sit on row 7
# (the Hepatica nobilis / Blåsippa record) and now belongs on row 5,
# where that record lives after the rotation.
$ws.Cells.Item(7, 12).ClearContents()
$ws.Cells.Item(5, 12).Value = ""
